# "feat: added fake data for beta dicts"
#
# The "Variables" sheet had a placeholder row (name/text/Name of the child)
# that gets removed; the following "country" row shifts up to take its
# place (row 5 -> row 4). The "Variables" sheet also becomes the active
# (selected) sheet/tab instead of "Categories".

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")

# Remove the placeholder "name"/"text" row (row 4) - the "country" row
# that followed it shifts up into row 4.
$wsVariables.Rows.Item(4).Delete() | Out-Null

# Make "Variables" the active sheet/tab.
$wsVariables.Activate() | Out-Null

# Select the (new) row 4, matching the saved selection state.
$wsVariables.Rows.Item(4).Select() | Out-Null
